$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "B" and "C" sub-rows within each year group (A col labels stay
# attached to their own data, rows 3/4, 7/8, 11/12, 15/16 trade places).
$pairs = @(3, 7, 11, 15)
foreach ($r1 in $pairs) {
    $r2 = $r1 + 1
    $row1 = $ws.Range("A$r1`:E$r1").Value2
    $row2 = $ws.Range("A$r2`:E$r2").Value2
    $ws.Range("A$r1`:E$r1").Value2 = $row2
    $ws.Range("A$r2`:E$r2").Value2 = $row1
}

# Drop the now-unwanted "产销率" / "销售量" (F/G) columns entirely.
$ws.Columns("F:G").Delete()
